# Update TPM-derived values in the LR-pairs sheet (Gdf2-Acvr2b)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Receptor average/total expression values change, which cascades into
# the derived specificity and edge-weight columns.
$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 1.841987434754
$ws.Range("R2").Value = 16.577886912786
$ws.Range("S2").Value = 0.3057455162066235
$ws.Range("T2").Value = 0.3057455162066235

# Row 3: only the derived-specificity columns shift (recomputed vs. new row 2)
$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("S3").Value = 0.2805555239151429
$ws.Range("T3").Value = 0.2805555239151429

# Row 4: only the derived-specificity columns shift (recomputed vs. new row 2)
$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("S4").Value = 0.4136989598782336
$ws.Range("T4").Value = 0.4136989598782336
